# "Added last minute updates"
#
# Targets the first paragraph of the document (the
# **ID__AFFARS_5315_topic_5__ID** placeholder paragraph):
#   - give it a paragraph border with 5-twip spacing on all four sides
#   - bump its left indent from 120 -> 225 twips (6pt -> 11.25pt)
#   - rewrite the placeholder text to **ID__AFFARS_5315_371_4__ID**
#   - drop the trailing run that held a single literal space
#
$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# w:pBdr with w:top/left/bottom/right each w:space="5", no line (so no
# w:sz/w:val/w:color attributes get emitted on the border elements).
$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromRight = 5

# w:ind w:left="120" -> w:left="225" (twips); LeftIndent is in points.
$p1.Format.LeftIndent = 11.25

# Collapse the paragraph's two runs (text run + trailing-space run) into
# a single run carrying the updated placeholder text. Excluding the
# paragraph mark keeps this a plain in-paragraph text replace.
$r = $p1.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "**ID__AFFARS_5315_371_4__ID**"
